$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Numero(s) caso" column header (J1) is renamed to "Numero de caso"
$ws.Range("J1").Value = "Numero de caso"

# New column L is added with header "Documentos analizados"
$ws.Range("L1").Value = "Documentos analizados"

# Give the new column L the same kind of explicit width the others have
$ws.Columns.Item(12).ColumnWidth = 19.6

# Move/update the view: scroll so column E is the left-most visible column,
# and select I10 as the active cell (as in the edited workbook)
$ws.Activate()
$win = $ws.Application.ActiveWindow
$win.ScrollColumn = 5
$ws.Range("I10").Select()
